# Update "想去人数" (expected attendance) figures on the "展览" sheet
# and the matching rows on the "全部类型" sheet (which mirrors the same
# events), reflecting freshly scraped output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 3098
$wsExhibit.Range("F4").Value = 145
$wsExhibit.Range("F5").Value = 112

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 3098
$wsAll.Range("F8").Value = 145
$wsAll.Range("F10").Value = 112
